# Update the cryptos list with latest price/volume data.
#
# The Price column (D) stores values that look numeric (e.g. "1.001",
# "303.94") but must remain plain text, matching how the workbook already
# stores them. Assigning a bare numeric-looking string makes Excel's COM
# layer coerce the cell to a real number (and mangle values like "1.0000"
# down to "1"), so we prefix with a leading apostrophe to force text entry,
# then reset the cell style to "Normal" so the implicit quote-prefix
# formatting doesn't leave a stray style on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-PriceText "D2" "23.735.59"
$ws.Range("E2").Value = "  +0.92%  "

# Row 3 - Ethereum
Set-PriceText "D3" "1.658.17"
$ws.Range("E3").Value = "  +1.11%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "1.0000"
$ws.Range("E4").Value = "  +0.15%  "

# Row 5 - USDC
Set-PriceText "D5" "1.001"
$ws.Range("E5").Value = "  +0.12%  "

# Row 6 - BNB
Set-PriceText "D6" "303.94"
$ws.Range("E6").Value = "  -0.09%  "

# Row 7 - XRP
Set-PriceText "D7" "0.3811"
$ws.Range("E7").Value = "  +0.53%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.16%  "

# Row 9 - OKB
Set-PriceText "D9" "51.23"
$ws.Range("E9").Value = "  -1.07%  "

# Row 10 - now Dogecoin (was Polygon)
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-PriceText "D10" "0.08213"
$ws.Range("E10").Value = "  +0.16%  "

# Row 11 - now Polygon (was Dogecoin)
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-PriceText "D11" "1.242"
$ws.Range("E11").Value = "  +0.38%  "

# Row 12 - BinanceUSD
Set-PriceText "D12" "0.9999"
$ws.Range("E12").Value = "  +0.12%  "

# Row 13 - Solana
$ws.Range("E13").Value = "  +0.42%  "

# Row 14 - Polkadot
Set-PriceText "D14" "6.529"
$ws.Range("E14").Value = "  +0.92%  "

# Row 15 - Chainlink
Set-PriceText "D15" "7.431"
$ws.Range("E15").Value = "  +0.63%  "

# Row 16 - ShibaInu
Set-PriceText "D16" "0.00001235"
$ws.Range("E16").Value = "  -0.72%  "

# Row 17 - WrappedEther
Set-PriceText "D17" "1.655.62"
$ws.Range("E17").Value = "  +0.80%  "

# Row 18 - Litecoin
Set-PriceText "D18" "97.57"
$ws.Range("E18").Value = "  +2.36%  "

# Row 19 - TRON
Set-PriceText "D19" "0.06998"
$ws.Range("E19").Value = "  +0.79%  "

# Row 20 - Uniswap
Set-PriceText "D20" "6.837"
$ws.Range("E20").Value = "  +3.64%  "

# Row 21 - Avalanche
Set-PriceText "D21" "17.69"
$ws.Range("E21").Value = "  +1.12%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.13%  "

# Row 24 - WrappedBTC
Set-PriceText "D24" "23.737.02"
$ws.Range("E24").Value = "  +0.95%  "

# Row 25 - Toncoin
Set-PriceText "D25" "2.517"
$ws.Range("E25").Value = "  -0.09%  "

# Row 26 - LidoDAOToken
Set-PriceText "D26" "3.053"
$ws.Range("E26").Value = "  -0.36%  "

# Row 27 - EthereumClassic
Set-PriceText "D27" "21.31"
$ws.Range("E27").Value = "  +0.71%  "

# Row 28 - Monero
Set-PriceText "D28" "153.41"
$ws.Range("E28").Value = "  +1.19%  "

# Row 29 - HuobiToken
Set-PriceText "D29" "5.212"
$ws.Range("E29").Value = "  -1.16%  "

# Row 30 - BitcoinCash
Set-PriceText "D30" "134.51"
$ws.Range("E30").Value = "  +0.89%  "

# Row 31 - WrappedliquidstakedEther2.0
Set-PriceText "D31" "1.840.17"
$ws.Range("E31").Value = "  +1.48%  "

# Row 32 - Filecoin
Set-PriceText "D32" "6.966"
$ws.Range("E32").Value = "  +4.81%  "

# Row 33 - WEMIXTOKEN
Set-PriceText "D33" "2.198"
$ws.Range("E33").Value = "  +0.40%  "

# Row 34 - ImmutableX
Set-PriceText "D34" "1.065"
$ws.Range("E34").Value = "  +0.90%  "

# Row 35 - FraxShare
Set-PriceText "D35" "11.82"
$ws.Range("E35").Value = "  +4.18%  "

# Row 36 - VeChain
Set-PriceText "D36" "0.02818"
$ws.Range("E36").Value = "  +2.02%  "

# Row 37 - Algorand
Set-PriceText "D37" "0.2530"
$ws.Range("E37").Value = "  +1.28%  "

# Row 38 - InternetComputer(DFINITY)
Set-PriceText "D38" "6.114"
$ws.Range("E38").Value = "  +1.47%  "

# Row 39 - Stellar
Set-PriceText "D39" "0.08786"
$ws.Range("E39").Value = "  +0.10%  "

# Row 40 - Hedera
Set-PriceText "D40" "0.07061"
$ws.Range("E40").Value = "  -0.96%  "

# Row 41 - Aptos
$ws.Range("E41").Value = "  +7.21%  "

# Row 42 - TheSandbox
Set-PriceText "D42" "0.7046"
$ws.Range("E42").Value = "  -0.72%  "

# Row 43 - TrustWalletToken
Set-PriceText "D43" "1.336"
$ws.Range("E43").Value = "  -0.57%  "

# Row 44 - EnergySwap
Set-PriceText "D44" "16.12"
$ws.Range("E44").Value = "  +2.52%  "

# Row 45 - Decentraland
Set-PriceText "D45" "0.6537"
$ws.Range("E45").Value = "  -0.27%  "

# Row 46 - NEARProtocol
Set-PriceText "D46" "2.325"
$ws.Range("E46").Value = "  +1.65%  "

# Row 47 - Frax
Set-PriceText "D47" "1.001"
$ws.Range("E47").Value = "  +0.14%  "

# Row 48 - PancakeSwap
Set-PriceText "D48" "3.981"
$ws.Range("E48").Value = "  +0.31%  "

# Row 49 - Cronos
Set-PriceText "D49" "0.07948"

# Row 50 - Quant
Set-PriceText "D50" "128.36"
$ws.Range("E50").Value = "  +0.54%  "

# Row 51 - Flow
Set-PriceText "D51" "1.186"
$ws.Range("E51").Value = "  -0.84%  "
